$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.150.66'
$ws.Range("E2").Value = '  -4.40%  '
$ws.Range("D3").Value = '1.655.23'
$ws.Range("E3").Value = '  -3.30%  '
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("D5").Value = '215.37'
$ws.Range("E5").Value = '  -4.11%  '
$ws.Range("D6").Value = '0.5089'
$ws.Range("E6").Value = '  -4.02%  '
$ws.Range("E7").Value = '  +0.16%  '
$ws.Range("E8").Value = '  -3.30%  '
$ws.Range("D9").Value = '0.06415'
$ws.Range("E9").Value = '  -4.22%  '
$ws.Range("E10").Value = '  -4.75%  '
$ws.Range("D11").Value = '0.07795'
$ws.Range("E11").Value = '  +1.46%  '
$ws.Range("D12").Value = '1.654.08'
$ws.Range("E12").Value = '  -3.17%  '
$ws.Range("D13").Value = '4.283'
$ws.Range("E13").Value = '  -5.17%  '
$ws.Range("D14").Value = '1.882.45'
$ws.Range("E14").Value = '  -3.33%  '
$ws.Range("E15").Value = '  -5.74%  '
$ws.Range("E16").Value = '  -2.81%  '
$ws.Range("D17").Value = '63.99'
$ws.Range("E17").Value = '  -6.08%  '
$ws.Range("D18").Value = '26.162.31'
$ws.Range("E18").Value = '  -4.39%  '
$ws.Range("E19").Value = '  +0.17%  '
$ws.Range("D20").Value = '208.67'
$ws.Range("E20").Value = '  -6.56%  '
$ws.Range("D21").Value = '4.408'
$ws.Range("E21").Value = '  -5.03%  '
$ws.Range("E22").Value = '  -3.49%  '
$ws.Range("D23").Value = '5.995'
$ws.Range("E23").Value = '  -0.20%  '
$ws.Range("E24").Value = '  +0.13%  '
$ws.Range("E25").Value = '  -0.53%  '
$ws.Range("D26").Value = '1.738'
$ws.Range("E26").Value = '  +2.67%  '
$ws.Range("D27").Value = '0.1172'
$ws.Range("E27").Value = '  -3.55%  '
$ws.Range("D28").Value = '6.967'
$ws.Range("E28").Value = '  -3.98%  '
$ws.Range("D29").Value = '15.80'
$ws.Range("E29").Value = '  -2.93%  '
$ws.Range("D30").Value = '0.05108'
$ws.Range("E30").Value = '  -4.96%  '
$ws.Range("E31").Value = '  -4.12%  '
$ws.Range("D32").Value = '3.351'
$ws.Range("E32").Value = '  -3.12%  '
$ws.Range("E33").Value = '  -6.51%  '
$ws.Range("D34").Value = '1.568'
$ws.Range("E34").Value = '  -4.53%  '
$ws.Range("D35").Value = '2.748'
$ws.Range("E35").Value = '  -4.12%  '
$ws.Range("D36").Value = '2.372'
$ws.Range("E36").Value = '  -0.85%  '
$ws.Range("E37").Value = '  -2.61%  '
$ws.Range("D38").Value = '0.5689'
$ws.Range("E38").Value = '  -3.07%  '
$ws.Range("D39").Value = '1.154.89'
$ws.Range("E39").Value = '  +5.81%  '
$ws.Range("D40").Value = '0.01589'
$ws.Range("E40").Value = '  -3.07%  '
$ws.Range("D41").Value = '2.561'
$ws.Range("E41").Value = '  -0.19%  '
$ws.Range("E42").Value = '  +0.09%  '
$ws.Range("D43").Value = '0.8337'
$ws.Range("E43").Value = '  -1.25%  '
$ws.Range("D44").Value = '5.636'
$ws.Range("E44").Value = '  -3.26%  '
$ws.Range("D45").Value = '100.26'
$ws.Range("E45").Value = '  -0.63%  '
$ws.Range("D46").Value = '1.792.25'
$ws.Range("E46").Value = '  -3.35%  '
$ws.Range("E47").Value = '  +1.55%  '
$ws.Range("D48").Value = '0.4548'
$ws.Range("E48").Value = '  +0.33%  '
$ws.Range("D49").Value = '55.66'
$ws.Range("E49").Value = '  -3.97%  '
$ws.Range("D50").Value = '1.005'
$ws.Range("E50").Value = '  +0.04%  '
$ws.Range("D51").Value = '7.860'
$ws.Range("E51").Value = '  -3.17%  '

Write-Host "Updated cryptos list"
